# feat: add 2022-Q1 data
#
# - Insert a new worksheet "2022-Q1" (fund-holdings detail) between the
#   existing "2021-Q4" sheet and the "总计" (total) summary sheet.
# - Add a "2022-Q1" summary row at the top of the "总计" sheet's data
#   (pushing the existing "2021-Q4" row down one row).

function Set-TextCell($ws, $addr, $val) {
    # Force text storage so numeric-looking strings (fund codes like
    # "006102", ratios like "48.75") keep their literal representation
    # instead of being auto-coerced to numbers (which would e.g. drop
    # leading zeros or normalise "48.75" -> 48.75 as a float).
    $ws.Range($addr).NumberFormat = "@"
    $ws.Range($addr).Value = $val
    $ws.Range($addr).Style = "Normal"
}

$wb = $excel.ActiveWorkbook

$q4Sheet = $wb.Worksheets.Item(1)
$totalSheetBefore = $wb.Worksheets.Item(2)

# ---------------------------------------------------------------------
# 1. Insert the new "2022-Q1" sheet right before the "总计" sheet, so the
#    final order is: 2021-Q4, 2022-Q1, 总计.
#    NOTE: inserting shifts worksheet positions, so the "总计" sheet must
#    be re-fetched by its *new* index afterwards -- the pre-insert
#    reference ends up tracking whichever sheet now sits at that old
#    position (i.e. the freshly inserted one), not the original "总计".
# ---------------------------------------------------------------------
$q1Sheet = $wb.Worksheets.Add($totalSheetBefore)
$q1Sheet.Name = "2022-Q1"

$totalSheet = $wb.Worksheets.Item(3)

# ---------------------------------------------------------------------
# 2. Populate the "2022-Q1" sheet with the fund-holdings table.
#
#    The header row + first data row's formatting (bold/centered/bordered
#    header in B1:H1, bordered index cell in column A) mirrors the
#    "2021-Q4" sheet exactly, so seed the new sheet by copying that
#    sheet's A1:H2 block (brings along both the styles and, for the
#    header text, the correct values too) and then overwrite the data
#    cells with the 2022-Q1 figures. Copying row 2's style down to rows
#    3-9 keeps every data row consistently formatted.
# ---------------------------------------------------------------------
$q4Sheet.Range("A1:H2").Copy($q1Sheet.Range("A1:H2"))
for ($i = 3; $i -le 9; $i++) {
    $q1Sheet.Range("A2:H2").Copy($q1Sheet.Range("A" + $i + ":H" + $i))
}

$rows = @(
    @{A=0; B="006102"; C="浙商丰利增强债券";                 D="48.75"; E="47.92"; F="2.19"; G="1.0676"; H=6},
    @{A=1; B="688888"; C="浙商聚潮产业成长混合";             D="8.25";  E="93.40"; F="5.11"; G="0.4216"; H=7},
    @{A=2; B="010381"; C="浙商智选价值混合A";                D="2.92";  E="93.43"; F="4.81"; G="0.1405"; H=9},
    @{A=3; B="010382"; C="浙商智选价值混合C";                D="0.34";  E="93.43"; F="4.81"; G="0.0164"; H=9},
    @{A=4; B="512590"; C="浦银安盛中证高股息精选ETF";        D="0.59";  E="96.43"; F="2.51"; G="0.0148"; H=4},
    @{A=5; B="006143"; C="恒生前海中证质量成长低波动指数A";  D="0.06";  E="94.34"; F="2.22"; G="0.0013"; H=10},
    @{A=6; B="005770"; C="信达澳银中证沪港深高股息精选指数"; D="0.01";  E="92.47"; F="2.16"; G="0.0002"; H=10},
    @{A=7; B="006144"; C="恒生前海中证质量成长低波动指数C";  D="0.01";  E="94.34"; F="2.22"; G="0.0002"; H=10}
)

$r = 2
foreach ($row in $rows) {
    $q1Sheet.Range("A$r").Value = $row.A

    Set-TextCell $q1Sheet "B$r" $row.B
    Set-TextCell $q1Sheet "C$r" $row.C
    Set-TextCell $q1Sheet "D$r" $row.D
    Set-TextCell $q1Sheet "E$r" $row.E
    Set-TextCell $q1Sheet "F$r" $row.F
    Set-TextCell $q1Sheet "G$r" $row.G

    $q1Sheet.Range("H$r").Value = $row.H

    $r = $r + 1
}

# ---------------------------------------------------------------------
# 3. Update the "总计" sheet: push the existing "2021-Q4" row from row 2
#    to row 3 (copy keeps its style), then write the new "2022-Q1"
#    summary into row 2.
# ---------------------------------------------------------------------
$totalSheet.Range("A2:D2").Copy($totalSheet.Range("A3:D3"))
$totalSheet.Range("A3").Value = 1

$totalSheet.Range("A2").Value = 0
$totalSheet.Range("B2").Value = "2022-Q1"
$totalSheet.Range("C2").Value = 8
$totalSheet.Range("D2").Value = 1.66

Write-Host "2022-Q1 data added"
